$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = '51.957.72'
$ws.Cells.Item(2, 5).Value2 = '  +0.27%  '

$ws.Cells.Item(3, 4).Value2 = '2.782.08'
$ws.Cells.Item(3, 5).Value2 = '  -1.15%  '

$ws.Cells.Item(4, 5).Value2 = '  -0.05%  '

$ws.Cells.Item(5, 4).Value2 = "'357.33"
$ws.Cells.Item(5, 4).Style = 'Normal'
$ws.Cells.Item(5, 5).Value2 = '  +0.99%  '

$ws.Cells.Item(6, 4).Value2 = "'109.86"
$ws.Cells.Item(6, 4).Style = 'Normal'
$ws.Cells.Item(6, 5).Value2 = '  -3.02%  '

$ws.Cells.Item(7, 4).Value2 = "'0.567"
$ws.Cells.Item(7, 4).Style = 'Normal'
$ws.Cells.Item(7, 5).Value2 = '  +2.75%  '

$ws.Cells.Item(8, 4).Value2 = "'0.999"
$ws.Cells.Item(8, 4).Style = 'Normal'
$ws.Cells.Item(8, 5).Value2 = '  -0.02%  '

$ws.Cells.Item(9, 5).Value2 = '  -0.86%  '

$ws.Cells.Item(10, 4).Value2 = "'40.15"
$ws.Cells.Item(10, 4).Style = 'Normal'
$ws.Cells.Item(10, 5).Value2 = '  -3.31%  '

$ws.Cells.Item(11, 4).Value2 = "'0.0853"
$ws.Cells.Item(11, 4).Style = 'Normal'
$ws.Cells.Item(11, 5).Value2 = '  +0.19%  '

$ws.Cells.Item(12, 5).Value2 = '  +0.72%  '

$ws.Cells.Item(13, 4).Value2 = "'19.40"
$ws.Cells.Item(13, 4).Style = 'Normal'
$ws.Cells.Item(13, 5).Value2 = '  -3.00%  '

$ws.Cells.Item(14, 4).Value2 = "'7.63"
$ws.Cells.Item(14, 4).Style = 'Normal'
$ws.Cells.Item(14, 5).Value2 = '  -0.95%  '

$ws.Cells.Item(15, 4).Value2 = '3.216.06'
$ws.Cells.Item(15, 5).Value2 = '  -0.66%  '

$ws.Cells.Item(16, 4).Value2 = '2.730.87'
$ws.Cells.Item(16, 5).Value2 = '  -3.62%  '

$ws.Cells.Item(17, 4).Value2 = "'0.930"
$ws.Cells.Item(17, 4).Style = 'Normal'
$ws.Cells.Item(17, 5).Value2 = '  +3.85%  '

$ws.Cells.Item(18, 4).Value2 = '51.878.61'
$ws.Cells.Item(18, 5).Value2 = '  +0.19%  '

$ws.Cells.Item(19, 5).Value2 = '  +0.05%  '

$ws.Cells.Item(20, 5).Value2 = '  +0.01%  '

$ws.Cells.Item(21, 4).Value2 = "'13.07"
$ws.Cells.Item(21, 4).Style = 'Normal'
$ws.Cells.Item(21, 5).Value2 = '  -3.46%  '

$ws.Cells.Item(22, 4).Value2 = '0.0₃0977'
$ws.Cells.Item(22, 5).Value2 = '  -1.65%  '

$ws.Cells.Item(23, 4).Value2 = "'274.23"
$ws.Cells.Item(23, 4).Style = 'Normal'
$ws.Cells.Item(23, 5).Value2 = '  +1.64%  '

$ws.Cells.Item(24, 5).Value2 = '  +0.22%  '

$ws.Cells.Item(25, 4).Value2 = "'2.74"
$ws.Cells.Item(25, 4).Style = 'Normal'
$ws.Cells.Item(25, 5).Value2 = '  -1.44%  '

$ws.Cells.Item(26, 4).Value2 = "'26.59"
$ws.Cells.Item(26, 4).Style = 'Normal'
$ws.Cells.Item(26, 5).Value2 = '  -0.32%  '

$ws.Cells.Item(27, 5).Value2 = '  -0.05%  '

$ws.Cells.Item(28, 4).Value2 = "'10.16"
$ws.Cells.Item(28, 4).Style = 'Normal'
$ws.Cells.Item(28, 5).Value2 = '  -1.30%  '

$ws.Cells.Item(29, 2).Value2 = 'Toncoin'
$ws.Cells.Item(29, 3).Value2 = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(29, 4).Value2 = "'2.22"
$ws.Cells.Item(29, 4).Style = 'Normal'
$ws.Cells.Item(29, 5).Value2 = '  -1.12%  '

$ws.Cells.Item(30, 2).Value2 = 'Kaspa'
$ws.Cells.Item(30, 3).Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(30, 4).Value2 = "'0.143"
$ws.Cells.Item(30, 4).Style = 'Normal'
$ws.Cells.Item(30, 5).Value2 = '  +2.43%  '

$ws.Cells.Item(31, 5).Value2 = '  +2.75%  '

$ws.Cells.Item(32, 4).Value2 = "'51.67"
$ws.Cells.Item(32, 4).Style = 'Normal'
$ws.Cells.Item(32, 5).Value2 = '  +2.08%  '

$ws.Cells.Item(33, 4).Value2 = "'33.95"
$ws.Cells.Item(33, 4).Style = 'Normal'
$ws.Cells.Item(33, 5).Value2 = '  +0.80%  '

$ws.Cells.Item(34, 4).Value2 = "'5.71"
$ws.Cells.Item(34, 4).Style = 'Normal'
$ws.Cells.Item(34, 5).Value2 = '  -1.75%  '

$ws.Cells.Item(35, 4).Value2 = "'0.0845"
$ws.Cells.Item(35, 4).Style = 'Normal'
$ws.Cells.Item(35, 5).Value2 = '  +1.65%  '

$ws.Cells.Item(36, 5).Value2 = '  +7.40%  '

$ws.Cells.Item(37, 5).Value2 = '  -0.05%  '

$ws.Cells.Item(38, 5).Value2 = '  +1.18%  '

$ws.Cells.Item(39, 4).Value2 = "'18.12"
$ws.Cells.Item(39, 4).Style = 'Normal'
$ws.Cells.Item(39, 5).Value2 = '  -0.72%  '

$ws.Cells.Item(40, 5).Value2 = '  -4.16%  '

$ws.Cells.Item(41, 4).Value2 = "'2.53"
$ws.Cells.Item(41, 4).Style = 'Normal'
$ws.Cells.Item(41, 5).Value2 = '  -1.08%  '

$ws.Cells.Item(42, 5).Value2 = '  -0.37%  '

$ws.Cells.Item(43, 5).Value2 = '  -2.94%  '

$ws.Cells.Item(44, 4).Value2 = "'121.11"
$ws.Cells.Item(44, 4).Style = 'Normal'
$ws.Cells.Item(44, 5).Value2 = '  -4.72%  '

$ws.Cells.Item(45, 4).Value2 = "'22.06"
$ws.Cells.Item(45, 4).Style = 'Normal'
$ws.Cells.Item(45, 5).Value2 = '  -6.83%  '

$ws.Cells.Item(46, 4).Value2 = '2.070.57'
$ws.Cells.Item(46, 5).Value2 = '  -0.41%  '

$ws.Cells.Item(47, 5).Value2 = '  -2.87%  '

$ws.Cells.Item(48, 5).Value2 = '  -2.32%  '

$ws.Cells.Item(49, 4).Value2 = "'5.70"
$ws.Cells.Item(49, 4).Style = 'Normal'
$ws.Cells.Item(49, 5).Value2 = '  +0.64%  '

$ws.Cells.Item(50, 5).Value2 = '  -2.24%  '

$ws.Cells.Item(51, 4).Value2 = "'8.97"
$ws.Cells.Item(51, 4).Style = 'Normal'
$ws.Cells.Item(51, 5).Value2 = '  +0.66%  '
